$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 2435.25
$ws.Range("J17").Value = 2435.25
$ws.Range("L17").Value = 7305.75
$ws.Range("N17").Value = -7641.75
$ws.Range("H64").Value = 5624.75
$ws.Range("J64").Value = 9999.333000000001
$ws.Range("L64").Value = 9999.333000000001
$ws.Range("N64").Value = -10495.333
$ws.Range("H67").Value = 5624.75
$ws.Range("J67").Value = 9999.333000000001
$ws.Range("L67").Value = 9999.333000000001
$ws.Range("N67").Value = -11715.333
$ws.Range("H69").Value = 333339000
$ws.Range("J69").Value = 333339000
$ws.Range("L69").Value = 1000017000
$ws.Range("N69").Value = -1000018748
$ws.Range("H72").Value = 333339000
$ws.Range("J72").Value = 333339000
$ws.Range("L72").Value = 3000051000
$ws.Range("N72").Value = -3000059736
$ws.Range("H98").Value = 10146.695
$ws.Range("I98").Value = 12899.647
$ws.Range("J98").Value = 2346.6667
$ws.Range("K98").Value = 12899.647
$ws.Range("L98").Value = 2346.6667
$ws.Range("M98").Value = -11401.647
$ws.Range("N98").Value = -5342.6667
$ws.Range("H122").Value = 10146.695
$ws.Range("I122").Value = 12899.647
$ws.Range("J122").Value = 2346.6667
$ws.Range("K122").Value = 38698.94100000001
$ws.Range("L122").Value = 7040.000100000001
$ws.Range("M122").Value = -36248.94100000001
$ws.Range("N122").Value = -11940.0001
$ws.Range("H129").Value = 1499.25
$ws.Range("J129").Value = 2998
$ws.Range("L129").Value = 8994
$ws.Range("N129").Value = -18994
$ws.Range("H138").Value = 6182017
$ws.Range("I138").Value = 10881.091
$ws.Range("J138").Value = 15879516
$ws.Range("K138").Value = 32643.273
$ws.Range("L138").Value = 47638548
$ws.Range("M138").Value = -27503.273
$ws.Range("N138").Value = -47648828
$ws.Range("H141").Value = 10716.357
$ws.Range("I141").Value = 10670.333
$ws.Range("K141").Value = 32010.999
$ws.Range("M141").Value = -26830.999

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4104.041
$ws.Range("I32").Value = 4002.0417
$ws.Range("J32").Value = 9000
$ws.Range("K32").Value = 4002.0417
$ws.Range("L32").Value = 9000
$ws.Range("M32").Value = -3715.0417
$ws.Range("N32").Value = -9574
$ws.Range("H61").Value = 8217.5625
$ws.Range("I61").Value = 8283.23
$ws.Range("J61").Value = 7933
$ws.Range("K61").Value = 8283.23
$ws.Range("L61").Value = 7933
$ws.Range("M61").Value = -8071.23
$ws.Range("N61").Value = -8357
$ws.Range("H74").Value = 3226.9
$ws.Range("I74").Value = 3283.75
$ws.Range("K74").Value = 3283.75
$ws.Range("M74").Value = -2409.75
$ws.Range("H77").Value = 3226.9
$ws.Range("I77").Value = 3283.75
$ws.Range("K77").Value = 16418.75
$ws.Range("M77").Value = -12050.75
$ws.Range("H122").Value = 4631656
$ws.Range("I122").Value = 5557387
$ws.Range("K122").Value = 16672161
$ws.Range("M122").Value = -16669711
$ws.Range("H132").Value = 3793.98
$ws.Range("I132").Value = 3667.3264
$ws.Range("K132").Value = 11001.9792
$ws.Range("M132").Value = -8471.9792
$ws.Range("H136").Value = 8217.5625
$ws.Range("I136").Value = 8283.23
$ws.Range("J136").Value = 7933
$ws.Range("K136").Value = 24849.69
$ws.Range("L136").Value = 23799
$ws.Range("M136").Value = -22299.69
$ws.Range("N136").Value = -28899

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 7474.421
$ws.Range("I20").Value = 10268.333
$ws.Range("J20").Value = 2684.8572
$ws.Range("K20").Value = 10268.333
$ws.Range("L20").Value = 2684.8572
$ws.Range("M20").Value = -10021.333
$ws.Range("N20").Value = -3178.8572
$ws.Range("H75").Value = 12074.75
$ws.Range("I75").Value = 12074.75
$ws.Range("K75").Value = 12074.75
$ws.Range("M75").Value = -11138.75
$ws.Range("H78").Value = 12074.75
$ws.Range("I78").Value = 12074.75
$ws.Range("K78").Value = 36224.25
$ws.Range("M78").Value = -31544.25
$ws.Range("H107").Value = 4078.1365
$ws.Range("I107").Value = 2362.9333
$ws.Range("K107").Value = 2362.9333
$ws.Range("M107").Value = -442.9333000000001
$ws.Range("H134").Value = 2685.2415
$ws.Range("I134").Value = 2619.6667
$ws.Range("J134").Value = 3000
$ws.Range("K134").Value = 7859.000100000001
$ws.Range("L134").Value = 9000
$ws.Range("M134").Value = -5324.000100000001
$ws.Range("N134").Value = -14070

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H6").Value = 1505.75
$ws.Range("I6").Value = 1547.1818
$ws.Range("K6").Value = 1547.1818
$ws.Range("M6").Value = -1434.1818
$ws.Range("H31").Value = 3756.5
$ws.Range("I31").Value = 2120.375
$ws.Range("J31").Value = 4691.4287
$ws.Range("K31").Value = 2120.375
$ws.Range("L31").Value = 4691.4287
$ws.Range("M31").Value = -1825.375
$ws.Range("N31").Value = -5281.4287
$ws.Range("H34").Value = 3756.5
$ws.Range("I34").Value = 2120.375
$ws.Range("J34").Value = 4691.4287
$ws.Range("K34").Value = 2120.375
$ws.Range("L34").Value = 4691.4287
$ws.Range("M34").Value = -1918.375
$ws.Range("N34").Value = -5095.4287
$ws.Range("H58").Value = 2637.0667
$ws.Range("I58").Value = 2539.7144
$ws.Range("K58").Value = 2539.7144
$ws.Range("M58").Value = -2336.7144
$ws.Range("H132").Value = 2045.7273
$ws.Range("I132").Value = 2033.7778
$ws.Range("K132").Value = 6101.3334
$ws.Range("M132").Value = -3571.3334
$ws.Range("H134").Value = 1904.931
$ws.Range("I134").Value = 1927.5927
$ws.Range("K134").Value = 5782.7781
$ws.Range("M134").Value = -3247.7781
$ws.Range("H136").Value = 2637.0667
$ws.Range("I136").Value = 2539.7144
$ws.Range("K136").Value = 7619.1432
$ws.Range("M136").Value = -5069.1432
$ws.Range("H141").Value = 519999.75
$ws.Range("J141").Value = 519999.75
$ws.Range("L141").Value = 519999.75
$ws.Range("N141").Value = -530359.75

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H75").Value = 2050.6667
$ws.Range("J75").Value = 2683.75
$ws.Range("L75").Value = 8051.25
$ws.Range("N75").Value = -10047.25
$ws.Range("H78").Value = 2050.6667
$ws.Range("J78").Value = 2683.75
$ws.Range("L78").Value = 24153.75
$ws.Range("N78").Value = -34137.75
$ws.Range("H80").Value = 2416.6667
$ws.Range("J80").Value = 2416.6667
$ws.Range("L80").Value = 7250.000100000001
$ws.Range("N80").Value = -9122.000100000001
$ws.Range("H83").Value = 2416.6667
$ws.Range("J83").Value = 2416.6667
$ws.Range("L83").Value = 21750.0003
$ws.Range("N83").Value = -31110.0003

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 4982.489
$ws.Range("I80").Value = 3798.6072
$ws.Range("J80").Value = 6932.4116
$ws.Range("K80").Value = 3798.6072
$ws.Range("L80").Value = 6932.4116
$ws.Range("M80").Value = -2800.6072
$ws.Range("N80").Value = -8928.411599999999
$ws.Range("H83").Value = 4982.489
$ws.Range("I83").Value = 3798.6072
$ws.Range("J83").Value = 6932.4116
$ws.Range("K83").Value = 18993.036
$ws.Range("L83").Value = 34662.058
$ws.Range("M83").Value = -14001.036
$ws.Range("N83").Value = -44646.058
$ws.Range("H126").Value = 2624.2856
$ws.Range("I126").Value = 1867.5
$ws.Range("J126").Value = 3633.3333
$ws.Range("K126").Value = 5602.5
$ws.Range("L126").Value = 10899.9999
$ws.Range("M126").Value = -3132.5
$ws.Range("N126").Value = -15839.9999
$ws.Range("H132").Value = 8305.362999999999
$ws.Range("I132").Value = 4545.125
$ws.Range("J132").Value = 18332.666
$ws.Range("K132").Value = 13635.375
$ws.Range("L132").Value = 54997.99800000001
$ws.Range("M132").Value = -11105.375
$ws.Range("N132").Value = -60057.99800000001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H53").Value = 22168.2
$ws.Range("I53").Value = 25615.334
$ws.Range("J53").Value = 16997.5
$ws.Range("K53").Value = 25615.334
$ws.Range("L53").Value = 16997.5
$ws.Range("M53").Value = -25097.334
$ws.Range("N53").Value = -18033.5
$ws.Range("H55").Value = 1112.6945
$ws.Range("I55").Value = 731.76
$ws.Range("J55").Value = 1978.4546
$ws.Range("K55").Value = 731.76
$ws.Range("L55").Value = 1978.4546
$ws.Range("M55").Value = -558.76
$ws.Range("N55").Value = -2324.4546
$ws.Range("H100").Value = 7725.1875
$ws.Range("I100").Value = 7245.727
$ws.Range("J100").Value = 8780
$ws.Range("K100").Value = 7245.727
$ws.Range("L100").Value = 8780
$ws.Range("M100").Value = -6704.727
$ws.Range("N100").Value = -9862
$ws.Range("H136").Value = 16134.775
$ws.Range("I136").Value = 3914.652
$ws.Range("K136").Value = 11743.956
$ws.Range("M136").Value = -9193.956

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H46").Value = 53692.25
$ws.Range("J46").Value = 53692.25
$ws.Range("L46").Value = 53692.25
$ws.Range("N46").Value = -54154.25
$ws.Range("H107").Value = 640.3333
$ws.Range("I107").Value = 589.5
$ws.Range("K107").Value = 1768.5
$ws.Range("M107").Value = 151.5
$ws.Range("H126").Value = 3583.4614
$ws.Range("I126").Value = 2899
$ws.Range("J126").Value = 5865
$ws.Range("K126").Value = 8697
$ws.Range("L126").Value = 17595
$ws.Range("M126").Value = -6227
$ws.Range("N126").Value = -22535
$ws.Range("H132").Value = 7019.064
$ws.Range("I132").Value = 7067.3486
$ws.Range("K132").Value = 21202.0458
$ws.Range("M132").Value = -18672.0458
$ws.Range("H134").Value = 53692.25
$ws.Range("J134").Value = 53692.25
$ws.Range("L134").Value = 161076.75
$ws.Range("N134").Value = -166146.75
$ws.Range("H136").Value = 6182.926
$ws.Range("I136").Value = 4581.1665
$ws.Range("K136").Value = 13743.4995
$ws.Range("M136").Value = -11193.4995
